$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.197.27"
$ws.Range("D3").Value = "1.914.27"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.71"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4064"
$ws.Range("E8").Value = "  +3.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08491"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.126"
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.75"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.37"
$ws.Range("E12").Value = "  +14.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.438"
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("D14").Value = "1.915.47"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.388"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.17"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001113"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.48"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.013"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "30.221.52"
$ws.Range("E23").Value = "  +5.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.227"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("D26").Value = "2.140.97"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.43"
$ws.Range("E27").Value = "  +4.24%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.28"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.409"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.81"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.101"
$ws.Range("E31").Value = "  +5.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1066"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.007"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.643"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02491"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06578"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2210"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.186"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.224"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.94"
$ws.Range("E40").Value = "  +7.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.806"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6519"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.238"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6141"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.39"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.739"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.076"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.27"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.164"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.50"
$ws.Range("E51").Value = "  +4.53%  "
